# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# commit ("Updated cryptos list on Mon Mar 25 23:55:52 UTC 2024 with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    # Force the cell to stay a text value even when $Value looks numeric
    # (e.g. "19.34"), matching the workbook convention where Price/Volume
    # columns are stored as plain text. Restore the default style afterwards
    # so no residual number-format is left on the cell.
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "69.934.11"
$ws.Range("E2").Value = "  +3.92%  "

# Row 3
$ws.Range("D3").Value = "3.589.06"
$ws.Range("E3").Value = "  +3.90%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
Set-TextValue $ws.Range("D5") "586.91"
$ws.Range("E5").Value = "  +3.26%  "

# Row 6
Set-TextValue $ws.Range("D6") "189.07"
$ws.Range("E6").Value = "  +2.89%  "

# Row 7
$ws.Range("E7").Value = "  +1.10%  "

# Row 8
$ws.Range("D8").Value = "3.582.57"
$ws.Range("E8").Value = "  +3.82%  "

# Row 9
$ws.Range("E9").Value = "  -0.05%  "

# Row 10
$ws.Range("E10").Value = "  -0.87%  "

# Row 11
$ws.Range("E11").Value = "  +1.53%  "

# Row 12
Set-TextValue $ws.Range("D12") "57.66"
$ws.Range("E12").Value = "  +3.59%  "

# Row 13
$ws.Range("E13").Value = "  +2.17%  "

# Row 14
Set-TextValue $ws.Range("D14") "9.73"
$ws.Range("E14").Value = "  +3.75%  "

# Row 15
$ws.Range("D15").Value = "4.159.03"
$ws.Range("E15").Value = "  +3.81%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "19.34"
$ws.Range("E16").Value = "  +4.32%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.588.21"
$ws.Range("E17").Value = "  +3.50%  "

# Row 18
$ws.Range("D18").Value = "69.845.75"
$ws.Range("E18").Value = "  +3.71%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.43"
$ws.Range("E19").Value = "  +3.31%  "

# Row 20
$ws.Range("E20").Value = "  +0.17%  "

# Row 21
$ws.Range("E21").Value = "  +3.22%  "

# Row 22
Set-TextValue $ws.Range("D22") "487.73"
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
Set-TextValue $ws.Range("D23") "17.49"
$ws.Range("E23").Value = "  +15.68%  "

# Row 24
$ws.Range("E24").Value = "  +7.97%  "

# Row 25
$ws.Range("E25").Value = "  +5.44%  "

# Row 26
Set-TextValue $ws.Range("D26") "90.30"
$ws.Range("E26").Value = "  +0.63%  "

# Row 27
Set-TextValue $ws.Range("D27") "3.09"
$ws.Range("E27").Value = "  +3.94%  "

# Row 28
Set-TextValue $ws.Range("D28") "11.03"
$ws.Range("E28").Value = "  +0.98%  "

# Row 29
$ws.Range("E29").Value = "  +4.86%  "

# Row 30
Set-TextValue $ws.Range("D30") "32.21"
$ws.Range("E30").Value = "  +1.88%  "

# Row 31
Set-TextValue $ws.Range("D31") "7.45"
$ws.Range("E31").Value = "  +7.05%  "

# Row 32
Set-TextValue $ws.Range("D32") "623.26"
$ws.Range("E32").Value = "  +3.75%  "

# Row 33
$ws.Range("E33").Value = "  +4.87%  "

# Row 34
$ws.Range("E34").Value = "  +6.43%  "

# Row 35
Set-TextValue $ws.Range("D35") "65.03"
$ws.Range("E35").Value = "  +3.11%  "

# Row 36
$ws.Range("E36").Value = "  +3.66%  "

# Row 37
Set-TextValue $ws.Range("D37") "1.00"
$ws.Range("E37").Value = "  +0.07%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.402"
$ws.Range("E38").Value = "  +3.49%  "

# Row 39
Set-TextValue $ws.Range("D39") "37.81"
$ws.Range("E39").Value = "  +3.55%  "

# Row 40
$ws.Range("E40").Value = "  -1.53%  "

# Row 41
Set-TextValue $ws.Range("D41") "3.62"
$ws.Range("E41").Value = "  -1.33%  "

# Row 42
$ws.Range("D42").Value = "3.293.11"
$ws.Range("E42").Value = "  +4.45%  "

# Row 43
$ws.Range("E43").Value = "  +4.69%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0446"
$ws.Range("E44").Value = "  +4.67%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.66"
$ws.Range("E45").Value = "  +2.88%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.34"
$ws.Range("E46").Value = "  +1.88%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.137"
$ws.Range("E47").Value = "  +1.12%  "

# Row 48
Set-TextValue $ws.Range("D48") "9.04"
$ws.Range("E48").Value = "  +3.38%  "

# Row 49
$ws.Range("E49").Value = "  +5.87%  "

# Row 50
Set-TextValue $ws.Range("D50") "2.69"
$ws.Range("E50").Value = "  -4.76%  "

# Row 51
$ws.Range("E51").Value = "  -0.08%  "

